# "Fixed typo in dates." — update the session date on the title slide
# from "Wednesday, February 19, 2013" to "Wednesday, February 19, 2014".
#
# In the canonical OOXML, the original paragraph was two runs:
#   run1: "Wednesday, February 19, "
#   run2: "2013"
# and the corrected paragraph is three runs (the comma/space became its
# own run) plus a trailing endParaRPr:
#   run1: "Wednesday, February 19"
#   run2: ", "
#   run3: "2014"
#
# We reproduce that by splitting the text range at the exact same
# boundaries (touching a Font property is enough to force a run split
# without altering formatting) and then updating the year text in place.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the shape whose text contains the date, rather than assuming a
# fixed shape index.
$needle = "Wednesday, February 19, 2013"
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame) {
        if ($candidate.TextFrame.TextRange.Text.IndexOf($needle) -ge 0) {
            $shp = $candidate
        }
    }
}

$tr = $shp.TextFrame.TextRange

$full = $tr.Text
$idx0 = $full.IndexOf($needle)
$startW = $idx0 + 1   # TextRange.Characters is 1-based

$prefixLen = "Wednesday, February 19".Length
$commaLen = ", ".Length
$yearLen = "2013".Length

$commaStart = $startW + $prefixLen
$yearStart = $commaStart + $commaLen

# Split "Wednesday, February 19, " into its own run ("Wednesday, February 19")
# plus a new ", " run, by nudging a character-format property on the
# comma/space sub-range (same value in, same value out — it only forces
# the run boundary to be created at this exact span).
$commaRange = $tr.Characters($commaStart, $commaLen)
$commaRange.Font.Bold = $commaRange.Font.Bold

# Update the year in place; run boundaries already match "2013"/"2014".
$yearRange = $tr.Characters($yearStart, $yearLen)
$yearRange.Text = "2014"
